$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2466
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 2466
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = ""
$ws.Cells.Item(17, 13).Value = 7398
$ws.Cells.Item(17, 14).Value = -7734

$ws.Cells.Item(55, 8).Value = 382
$ws.Cells.Item(55, 10).Value = 544
$ws.Cells.Item(55, 12).Value = 544
$ws.Cells.Item(55, 14).Value = -972

$ws.Cells.Item(80, 8).Value = 871.871
$ws.Cells.Item(80, 9).Value = 899.9
$ws.Cells.Item(80, 10).Value = 858.5238000000001
$ws.Cells.Item(80, 11).Value = 2699.7
$ws.Cells.Item(80, 12).Value = 2575.5714
$ws.Cells.Item(80, 13).Value = -1701.7
$ws.Cells.Item(80, 14).Value = -4571.571400000001

$ws.Cells.Item(83, 8).Value = 871.871
$ws.Cells.Item(83, 9).Value = 899.9
$ws.Cells.Item(83, 10).Value = 858.5238000000001
$ws.Cells.Item(83, 11).Value = 8099.099999999999
$ws.Cells.Item(83, 12).Value = 7726.7142
$ws.Cells.Item(83, 13).Value = -3107.099999999999
$ws.Cells.Item(83, 14).Value = -17710.7142

$ws.Cells.Item(86, 8).Value = 2751.8076
$ws.Cells.Item(86, 10).Value = 1915.4286
$ws.Cells.Item(86, 12).Value = 1915.4286
$ws.Cells.Item(86, 14).Value = -4161.4286

$ws.Cells.Item(89, 8).Value = 2751.8076
$ws.Cells.Item(89, 10).Value = 1915.4286
$ws.Cells.Item(89, 12).Value = 9577.143
$ws.Cells.Item(89, 14).Value = -20809.143

$ws.Cells.Item(100, 8).Value = 3716.739
$ws.Cells.Item(100, 9).Value = 1299.0667
$ws.Cells.Item(100, 10).Value = 8249.875
$ws.Cells.Item(100, 11).Value = 1299.0667
$ws.Cells.Item(100, 12).Value = 8249.875
$ws.Cells.Item(100, 13).Value = -758.0667000000001
$ws.Cells.Item(100, 14).Value = -9331.875

$ws.Cells.Item(129, 8).Value = 1708.9231
$ws.Cells.Item(129, 10).Value = 2975.5
$ws.Cells.Item(129, 12).Value = 8926.5
$ws.Cells.Item(129, 14).Value = -18926.5

$ws.Cells.Item(137, 8).Value = 1297.4688
$ws.Cells.Item(137, 9).Value = 1140.6
$ws.Cells.Item(137, 11).Value = 3421.8
$ws.Cells.Item(137, 13).Value = -871.7999999999997

$ws.Cells.Item(138, 8).Value = 2964.3035
$ws.Cells.Item(138, 9).Value = 2002.4667
$ws.Cells.Item(138, 10).Value = 3316.195
$ws.Cells.Item(138, 11).Value = 6007.4001
$ws.Cells.Item(138, 12).Value = 9948.585000000001
$ws.Cells.Item(138, 13).Value = -867.4000999999998
$ws.Cells.Item(138, 14).Value = -20228.585

$ws.Cells.Item(141, 8).Value = 6349
$ws.Cells.Item(141, 9).Value = 6788.6
$ws.Cells.Item(141, 11).Value = 20365.8
$ws.Cells.Item(141, 13).Value = -15185.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(31, 8).Value = 3464.3333
$ws.Cells.Item(31, 9).Value = 3464.3333
$ws.Cells.Item(31, 11).Value = 3464.3333
$ws.Cells.Item(31, 13).Value = -3170.3333

$ws.Cells.Item(34, 8).Value = 127440.14
$ws.Cells.Item(34, 9).Value = 78675
$ws.Cells.Item(34, 10).Value = 164014
$ws.Cells.Item(34, 11).Value = 78675
$ws.Cells.Item(34, 12).Value = 164014
$ws.Cells.Item(34, 13).Value = -78404
$ws.Cells.Item(34, 14).Value = -164556

$ws.Cells.Item(82, 8).Value = 56090.5
$ws.Cells.Item(82, 10).Value = 56090.5
$ws.Cells.Item(82, 12).Value = 56090.5
$ws.Cells.Item(82, 14).Value = -56812.5

$ws.Cells.Item(85, 8).Value = 56090.5
$ws.Cells.Item(85, 10).Value = 56090.5
$ws.Cells.Item(85, 12).Value = 56090.5
$ws.Cells.Item(85, 14).Value = -58586.5

$ws.Cells.Item(122, 8).Value = 2593.7273
$ws.Cells.Item(122, 9).Value = 2518.4
$ws.Cells.Item(122, 11).Value = 7555.200000000001
$ws.Cells.Item(122, 13).Value = -5105.200000000001

$ws.Cells.Item(132, 8).Value = 1793.7878
$ws.Cells.Item(132, 10).Value = 1171
$ws.Cells.Item(132, 12).Value = 3513
$ws.Cells.Item(132, 14).Value = -8573

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(46, 8).Value = 42214.332
$ws.Cells.Item(46, 10).Value = 42214.332
$ws.Cells.Item(46, 12).Value = 42214.332
$ws.Cells.Item(46, 14).Value = -42810.332

$ws.Cells.Item(86, 8).Value = 1189273.5
$ws.Cells.Item(86, 9).Value = 1853.3334
$ws.Cells.Item(86, 11).Value = 1853.3334
$ws.Cells.Item(86, 13).Value = -730.3334

$ws.Cells.Item(89, 8).Value = 1189273.5
$ws.Cells.Item(89, 9).Value = 1853.3334
$ws.Cells.Item(89, 11).Value = 9266.666999999999
$ws.Cells.Item(89, 13).Value = -3650.666999999999

$ws.Cells.Item(102, 8).Value = 7367.5
$ws.Cells.Item(102, 9).Value = 7367.5
$ws.Cells.Item(102, 11).Value = 7367.5
$ws.Cells.Item(102, 13).Value = -4122.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 940.9583
$ws.Cells.Item(22, 9).Value = 649.6
$ws.Cells.Item(22, 11).Value = 649.6
$ws.Cells.Item(22, 13).Value = -299.6

$ws.Cells.Item(51, 8).Value = 59510
$ws.Cells.Item(51, 9).Value = 66013.336
$ws.Cells.Item(51, 11).Value = 66013.336
$ws.Cells.Item(51, 13).Value = -65277.336

$ws.Cells.Item(58, 8).Value = 2282.7222
$ws.Cells.Item(58, 9).Value = 2270
$ws.Cells.Item(58, 11).Value = 2270
$ws.Cells.Item(58, 13).Value = -2067

$ws.Cells.Item(61, 8).Value = 59510
$ws.Cells.Item(61, 9).Value = 66013.336
$ws.Cells.Item(61, 11).Value = 66013.336
$ws.Cells.Item(61, 13).Value = -65665.336

$ws.Cells.Item(99, 8).Value = 4539.8
$ws.Cells.Item(99, 9).Value = 3999.6667
$ws.Cells.Item(99, 10).Value = 5350
$ws.Cells.Item(99, 11).Value = 3999.6667
$ws.Cells.Item(99, 12).Value = 5350
$ws.Cells.Item(99, 13).Value = -2501.6667
$ws.Cells.Item(99, 14).Value = -8346

$ws.Cells.Item(105, 8).Value = 4604.5
$ws.Cells.Item(105, 9).Value = 4877.6
$ws.Cells.Item(105, 11).Value = 4877.6
$ws.Cells.Item(105, 13).Value = -3130.6

$ws.Cells.Item(126, 8).Value = 4539.8
$ws.Cells.Item(126, 9).Value = 3999.6667
$ws.Cells.Item(126, 10).Value = 5350
$ws.Cells.Item(126, 11).Value = 11999.0001
$ws.Cells.Item(126, 12).Value = 16050
$ws.Cells.Item(126, 13).Value = -9529.000100000001
$ws.Cells.Item(126, 14).Value = -20990

$ws.Cells.Item(136, 8).Value = 2282.7222
$ws.Cells.Item(136, 9).Value = 2270
$ws.Cells.Item(136, 11).Value = 6810
$ws.Cells.Item(136, 13).Value = -4260

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 2236
$ws.Cells.Item(18, 9).Value = 823.8
$ws.Cells.Item(18, 11).Value = 2471.4
$ws.Cells.Item(18, 13).Value = -2302.4

$ws.Cells.Item(29, 8).Value = 48.142857
$ws.Cells.Item(29, 9).Value = 50.5
$ws.Cells.Item(29, 10).Value = 47.2
$ws.Cells.Item(29, 11).Value = 151.5
$ws.Cells.Item(29, 12).Value = 141.6
$ws.Cells.Item(29, 13).Value = 125.5
$ws.Cells.Item(29, 14).Value = -695.6

$ws.Cells.Item(140, 8).Value = 1212.25
$ws.Cells.Item(140, 9).Value = 1212.25
$ws.Cells.Item(140, 11).Value = 3636.75
$ws.Cells.Item(140, 13).Value = 1543.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 5692.1763
$ws.Cells.Item(80, 9).Value = 4412.6665
$ws.Cells.Item(80, 11).Value = 4412.6665
$ws.Cells.Item(80, 13).Value = -3414.6665

$ws.Cells.Item(83, 8).Value = 5692.1763
$ws.Cells.Item(83, 9).Value = 4412.6665
$ws.Cells.Item(83, 11).Value = 22063.3325
$ws.Cells.Item(83, 13).Value = -17071.3325

$ws.Cells.Item(132, 8).Value = 2179.762
$ws.Cells.Item(132, 9).Value = 1903.3529
$ws.Cells.Item(132, 10).Value = 3354.5
$ws.Cells.Item(132, 11).Value = 5710.0587
$ws.Cells.Item(132, 12).Value = 10063.5
$ws.Cells.Item(132, 13).Value = -3180.0587
$ws.Cells.Item(132, 14).Value = -15123.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3194.8667
$ws.Cells.Item(22, 9).Value = 2765.2856
$ws.Cells.Item(22, 11).Value = 2765.2856
$ws.Cells.Item(22, 13).Value = -2470.2856

$ws.Cells.Item(27, 8).Value = 3194.8667
$ws.Cells.Item(27, 9).Value = 2765.2856
$ws.Cells.Item(27, 11).Value = 2765.2856
$ws.Cells.Item(27, 13).Value = -2658.2856

$ws.Cells.Item(93, 8).Value = 4154.115
$ws.Cells.Item(93, 9).Value = 949.2222
$ws.Cells.Item(93, 11).Value = 949.2222
$ws.Cells.Item(93, 14).Value = 298.7778

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(40, 8).Value = 25000
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 13).Value = ""

$ws.Cells.Item(68, 8).Value = 16135.5
$ws.Cells.Item(68, 10).Value = 16135.5
$ws.Cells.Item(68, 12).Value = 16135.5
$ws.Cells.Item(68, 14).Value = -17757.5

$ws.Cells.Item(71, 8).Value = 16135.5
$ws.Cells.Item(71, 10).Value = 16135.5
$ws.Cells.Item(71, 12).Value = 48406.5
$ws.Cells.Item(71, 14).Value = -56518.5

$ws.Cells.Item(96, 8).Value = 3016.5386
$ws.Cells.Item(96, 9).Value = 2197.4
$ws.Cells.Item(96, 10).Value = 3528.5
$ws.Cells.Item(96, 11).Value = 2197.4
$ws.Cells.Item(96, 12).Value = 3528.5
$ws.Cells.Item(96, 13).Value = -824.4000000000001
$ws.Cells.Item(96, 14).Value = -6274.5

$ws.Cells.Item(100, 8).Value = 1443.8695
$ws.Cells.Item(100, 9).Value = 1248
$ws.Cells.Item(100, 11).Value = 2496
$ws.Cells.Item(100, 13).Value = -1955

$ws.Cells.Item(107, 8).Value = 2325.5
$ws.Cells.Item(107, 9).Value = 2287.8333
$ws.Cells.Item(107, 10).Value = 2495
$ws.Cells.Item(107, 11).Value = 6863.499899999999
$ws.Cells.Item(107, 12).Value = 7485
$ws.Cells.Item(107, 13).Value = -4943.499899999999
$ws.Cells.Item(107, 14).Value = -11325

$ws.Cells.Item(113, 8).Value = 1347.619
$ws.Cells.Item(113, 9).Value = 1058.5
$ws.Cells.Item(113, 11).Value = 3175.5
$ws.Cells.Item(113, 13).Value = -1005.5

$ws.Cells.Item(122, 8).Value = 4330.5
$ws.Cells.Item(122, 9).Value = 2386.182
$ws.Cells.Item(122, 10).Value = 6706.8887
$ws.Cells.Item(122, 11).Value = 7158.545999999999
$ws.Cells.Item(122, 12).Value = 20120.6661
$ws.Cells.Item(122, 13).Value = -4708.545999999999
$ws.Cells.Item(122, 14).Value = -25020.6661

$ws.Cells.Item(132, 8).Value = 2263.2886
$ws.Cells.Item(132, 9).Value = 2060.2273
$ws.Cells.Item(132, 11).Value = 6180.6819
$ws.Cells.Item(132, 13).Value = -3650.6819
